$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Misc"
$ws.Range("B11").Value = "checkMagCondition"
$ws.Range("C11").Value = "A bool"
$ws.Range("D11").Value = "Whether the condition of the mag should be checked"
$ws.Range("E11").Value = "/"
$ws.Range("F11").Value = "/"
$ws.Range("G11").Value = $true

$ws.Range("A11").Borders.Item(7).LineStyle = 1
$ws.Range("A11").Borders.Item(8).LineStyle = 1
$ws.Range("A11").Borders.Item(9).LineStyle = 1

$ws.Range("B11:F11").Borders.Item(8).LineStyle = 1
$ws.Range("B11:F11").Borders.Item(9).LineStyle = 1

$ws.Range("G11").Borders.Item(8).LineStyle = 1
$ws.Range("G11").Borders.Item(9).LineStyle = 1
$ws.Range("G11").Borders.Item(10).LineStyle = 1
